$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cep (I2) and telefone (E2) cells to be text values instead of numbers
# (order matters for shared string table indexing to match target workbook)
$ws.Range("I2").Value = "29827-17"
$ws.Range("E2").Value = "11 12121212"

# Update the active selection to G8 (matches new sheetView selection)
$ws.Range("G8").Select()
